# Generate Report for Handback
# Replaces the two handback file identifiers (66a797f0... -> 8e53a332...,
# 6eb66b5f... -> ffff289983fe...) and refreshes the associated
# handoff/handback timestamps and xlf hash filenames across all three
# worksheets (Overview, zh-cn, de-de), keeping hyperlink targets the same
# but updating their displayed text to match the new file names.

$wb = $excel.ActiveWorkbook

$oldId1 = "66a797f0-fd63-4f90-8a46-51d289cb68d1"
$newId1 = "8e53a332-4536-4ddf-a71c-e42674a0f667"
$oldId2 = "6eb66b5f-02fa-45dd-8b8d-0b9b8b9d097f"
$newId2 = "ffff289983fe-742f-4d47-b7a3-bb177cf2ed72"

$newZhCnXlf = "$newId1.298952e52ae56235bff20abe86efc1d52178f8c3.zh-cn.xlf"
$newDeDeXlf = "$newId1.298952e52ae56235bff20abe86efc1d52178f8c3.de-de.xlf"

# ---------------------------------------------------------------------
# Sheet: Overview
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "$newId1.md"
$ov.Range("B2").Value = "e2e\$newId1.md"
$ov.Range("G2").Value = "2016-08-30 13:07:44"
$ov.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ov.Range("A3").Value = "$newId2.md"
$ov.Range("B3").Value = "e2e\$newId2.md"
$ov.Range("G3").Value = "2016-08-30 13:07:44"
$ov.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

foreach ($h in $ov.Hyperlinks) {
    $r = $h.Range.Row
    $c = $h.Range.Column
    if ($r -eq 2 -and $c -eq 2) {
        $h.TextToDisplay = "e2e\$newId1.md"
    } elseif ($r -eq 3 -and $c -eq 2) {
        $h.TextToDisplay = "e2e\$newId2.md"
    }
}

# ---------------------------------------------------------------------
# Sheet: zh-cn
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "$newId1.md"
$zh.Range("G2").Value = $newZhCnXlf
$zh.Range("H2").Value = "2016-08-30 13:07:38"
$zh.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("I2").Value = "$newId1.md"
$zh.Range("J2").Value = $newZhCnXlf
$zh.Range("K2").Value = "2016-08-30 13:07:55"
$zh.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$zh.Range("A3").Value = "$newId2.md"
$zh.Range("G3").Value = $newZhCnXlf
$zh.Range("H3").Value = "2016-08-30 13:07:38"
$zh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("I3").Value = "$newId2.md"
$zh.Range("J3").Value = $newZhCnXlf
$zh.Range("K3").Value = "2016-08-30 13:07:55"
$zh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

foreach ($h in $zh.Hyperlinks) {
    $r = $h.Range.Row
    $c = $h.Range.Column
    if ($r -eq 2 -and $c -eq 1) {
        $h.TextToDisplay = "$newId1.md"
    } elseif ($r -eq 2 -and $c -eq 9) {
        $h.TextToDisplay = "$newId1.md"
    } elseif ($r -eq 3 -and $c -eq 1) {
        $h.TextToDisplay = "$newId2.md"
    } elseif ($r -eq 3 -and $c -eq 9) {
        $h.TextToDisplay = "$newId2.md"
    }
}

# ---------------------------------------------------------------------
# Sheet: de-de
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "$newId1.md"
$de.Range("G2").Value = $newDeDeXlf
$de.Range("H2").Value = "2016-08-30 13:07:44"
$de.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("I2").Value = "$newId1.md"
$de.Range("J2").Value = $newDeDeXlf
$de.Range("K2").Value = "2016-08-30 13:08:11"
$de.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$de.Range("A3").Value = "$newId2.md"
$de.Range("G3").Value = $newDeDeXlf
$de.Range("H3").Value = "2016-08-30 13:07:44"
$de.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("I3").Value = "$newId2.md"
$de.Range("J3").Value = $newDeDeXlf
$de.Range("K3").Value = "2016-08-30 13:08:11"
$de.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

foreach ($h in $de.Hyperlinks) {
    $r = $h.Range.Row
    $c = $h.Range.Column
    if ($r -eq 2 -and $c -eq 1) {
        $h.TextToDisplay = "$newId1.md"
    } elseif ($r -eq 2 -and $c -eq 9) {
        $h.TextToDisplay = "$newId1.md"
    } elseif ($r -eq 3 -and $c -eq 1) {
        $h.TextToDisplay = "$newId2.md"
    } elseif ($r -eq 3 -and $c -eq 9) {
        $h.TextToDisplay = "$newId2.md"
    }
}
